$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 19:16"

# Libia moves up in the shared-string/table order: insert its row of data at row 116,
# shifting Hong Kong / Mali / Cuba down one row each (Libia's old row 119 entry is dropped).
# Row 4
$ws.Range("B4").Value = 4278067
$ws.Range("C4").Value = 29740
$ws.Range("D4").Value = 2036752
$ws.Range("E4").Value = 2092348
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 477
$ws.Range("H4").Value = 148967

# Row 6
$ws.Range("B6").Value = 1383854
$ws.Range("C6").Value = 46832
$ws.Range("D6").Value = 884659
$ws.Range("E6").Value = 467106
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 683
$ws.Range("H6").Value = 32089

# Row 19
$ws.Range("B19").Value = 225173
$ws.Range("C19").Value = 921
$ws.Range("D19").Value = 208477
$ws.Range("E19").Value = 11100
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = 5596

# Row 21
$ws.Range("B21").Value = 206182
$ws.Range("C21").Value = 222
$ws.Range("D21").Value = 190400
$ws.Range("E21").Value = 6581
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 9201

# Row 24
$ws.Range("B24").Value = 113502
$ws.Range("C24").Value = 296
$ws.Range("D24").Value = 99098
$ws.Range("E24").Value = 5519
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 8885

# Row 31
$ws.Range("B31").Value = 80036
$ws.Range("C31").Value = 987
$ws.Range("D31").Value = 34544
$ws.Range("E31").Value = 39985
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 39
$ws.Range("H31").Value = 5507

# Row 60
$ws.Range("B60").Value = 26764
$ws.Range("C60").Value = 605
$ws.Range("D60").Value = 18076
$ws.Range("E60").Value = 7542
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 10
$ws.Range("H60").Value = 1146

# Row 61
$ws.Range("B61").Value = 25869
$ws.Range("C61").Value = 24
$ws.Range("D61").Value = 23364
$ws.Range("E61").Value = 741
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 1764

# Row 69
$ws.Range("B69").Value = 16643
$ws.Range("C69").Value = 375
$ws.Range("D69").Value = 7574
$ws.Range("E69").Value = 8791
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 4
$ws.Range("H69").Value = 278

# Row 116
$ws.Range("A116").Value = "Libia"
$ws.Range("B116").Value = 2547
$ws.Range("C116").Value = 123
$ws.Range("D116").Value = 510
$ws.Range("E116").Value = 1979
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 58

# Row 117
$ws.Range("A117").Value = "Hong Kong"
$ws.Range("B117").Value = 2506
$ws.Range("C117").Value = 133
$ws.Range("D117").Value = 1455
$ws.Range("E117").Value = 1033
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 2
$ws.Range("H117").Value = 18

# Row 118
$ws.Range("A118").Value = "Mali"
$ws.Range("B118").Value = 2503
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 1907
$ws.Range("E118").Value = 473
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 123

# Row 119
$ws.Range("A119").Value = "Cuba"
$ws.Range("B119").Value = 2469
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 2341
$ws.Range("E119").Value = 41
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 87

# Row 137
$ws.Range("B137").Value = 1443
$ws.Range("C137").Value = 18
$ws.Range("D137").Value = 1133
$ws.Range("E137").Value = 260
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 50

# Row 142
$ws.Range("B142").Value = 1154
$ws.Range("C142").Value = 8
$ws.Range("D142").Value = 1036
$ws.Range("E142").Value = 107
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 11

# Row 159
$ws.Range("B159").Value = 627
$ws.Range("C159").Value = 19
$ws.Range("D159").Value = 191
$ws.Range("E159").Value = 400
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 36

